$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "4+0=4"
$t.Cell(1,2).Range.Text = "23-4=19"
$t.Cell(1,3).Range.Text = "55+9=64"
$t.Cell(1,4).Range.Text = "11-11=0"
$t.Cell(1,5).Range.Text = "53+39=92"
$t.Cell(2,1).Range.Text = "73+24=97"
$t.Cell(2,2).Range.Text = "56+26=82"
$t.Cell(2,3).Range.Text = "13-3=10"
$t.Cell(2,4).Range.Text = "29+54=83"
$t.Cell(2,5).Range.Text = "83-76=7"
$t.Cell(3,1).Range.Text = "76-20=56"
$t.Cell(3,2).Range.Text = "67-50=17"
$t.Cell(3,3).Range.Text = "22+47=69"
$t.Cell(3,4).Range.Text = "44+29=73"
$t.Cell(3,5).Range.Text = "79+2=81"
$t.Cell(4,1).Range.Text = "46-31=15"
$t.Cell(4,2).Range.Text = "66-16=50"
$t.Cell(4,3).Range.Text = "31+53=84"
$t.Cell(4,4).Range.Text = "36-28=8"
$t.Cell(4,5).Range.Text = "42+52=94"
$t.Cell(5,1).Range.Text = "4+65=69"
$t.Cell(5,2).Range.Text = "65-11=54"
$t.Cell(5,3).Range.Text = "27+62=89"
$t.Cell(5,4).Range.Text = "62-33=29"
$t.Cell(5,5).Range.Text = "31+50=81"
$t.Cell(6,1).Range.Text = "75-21=54"
$t.Cell(6,2).Range.Text = "0+39=39"
$t.Cell(6,3).Range.Text = "52+33=85"
$t.Cell(6,4).Range.Text = "90-69=21"
$t.Cell(6,5).Range.Text = "23-10=13"
$t.Cell(7,1).Range.Text = "30+27=57"
$t.Cell(7,2).Range.Text = "25-18=7"
$t.Cell(7,3).Range.Text = "14+72=86"
$t.Cell(7,4).Range.Text = "83-5=78"
$t.Cell(7,5).Range.Text = "53-24=29"
$t.Cell(8,1).Range.Text = "42-12=30"
$t.Cell(8,2).Range.Text = "47+0=47"
$t.Cell(8,3).Range.Text = "64-36=28"
$t.Cell(8,4).Range.Text = "6+45=51"
$t.Cell(8,5).Range.Text = "99-17=82"
$t.Cell(9,1).Range.Text = "48-2=46"
$t.Cell(9,2).Range.Text = "1+65=66"
$t.Cell(9,3).Range.Text = "98-88=10"
$t.Cell(9,4).Range.Text = "51-23=28"
$t.Cell(9,5).Range.Text = "38+15=53"
$t.Cell(10,1).Range.Text = "34+21=55"
$t.Cell(10,2).Range.Text = "77-51=26"
$t.Cell(10,3).Range.Text = "65-65=0"
$t.Cell(10,4).Range.Text = "49+47=96"
$t.Cell(10,5).Range.Text = "94-39=55"
$t.Cell(11,1).Range.Text = "28+4=32"
$t.Cell(11,2).Range.Text = "45+7=52"
$t.Cell(11,3).Range.Text = "76+16=92"
$t.Cell(11,4).Range.Text = "55-23=32"
$t.Cell(11,5).Range.Text = "64+6=70"
$t.Cell(12,1).Range.Text = "61+1=62"
$t.Cell(12,2).Range.Text = "23+25=48"
$t.Cell(12,3).Range.Text = "69+19=88"
$t.Cell(12,4).Range.Text = "88-32=56"
$t.Cell(12,5).Range.Text = "51-29=22"
$t.Cell(13,1).Range.Text = "55+9=64"
$t.Cell(13,2).Range.Text = "35+41=76"
$t.Cell(13,3).Range.Text = "6+43=49"
$t.Cell(13,4).Range.Text = "30+65=95"
$t.Cell(13,5).Range.Text = "22+9=31"
$t.Cell(14,1).Range.Text = "45-33=12"
$t.Cell(14,2).Range.Text = "98-44=54"
$t.Cell(14,3).Range.Text = "23+17=40"
$t.Cell(14,4).Range.Text = "70+0=70"
$t.Cell(14,5).Range.Text = "89-35=54"
$t.Cell(15,1).Range.Text = "43-42=1"
$t.Cell(15,2).Range.Text = "75-45=30"
$t.Cell(15,3).Range.Text = "2+53=55"
$t.Cell(15,4).Range.Text = "42-18=24"
$t.Cell(15,5).Range.Text = "52+32=84"
$t.Cell(16,1).Range.Text = "42+42=84"
$t.Cell(16,2).Range.Text = "45+0=45"
$t.Cell(16,3).Range.Text = "67+11=78"
$t.Cell(16,4).Range.Text = "54-46=8"
$t.Cell(16,5).Range.Text = "53-32=21"
$t.Cell(17,1).Range.Text = "22+64=86"
$t.Cell(17,2).Range.Text = "56-30=26"
$t.Cell(17,3).Range.Text = "46-19=27"
$t.Cell(17,4).Range.Text = "92-78=14"
$t.Cell(17,5).Range.Text = "22+53=75"
$t.Cell(18,1).Range.Text = "52+44=96"
$t.Cell(18,2).Range.Text = "73-4=69"
$t.Cell(18,3).Range.Text = "56-20=36"
$t.Cell(18,4).Range.Text = "36+2=38"
$t.Cell(18,5).Range.Text = "36+57=93"
$t.Cell(19,1).Range.Text = "20+28=48"
$t.Cell(19,2).Range.Text = "4+81=85"
$t.Cell(19,3).Range.Text = "44-31=13"
$t.Cell(19,4).Range.Text = "36+41=77"
$t.Cell(19,5).Range.Text = "64+2=66"
$t.Cell(20,1).Range.Text = "21-9=12"
$t.Cell(20,2).Range.Text = "33+19=52"
$t.Cell(20,3).Range.Text = "34+31=65"
$t.Cell(20,4).Range.Text = "65-64=1"
$t.Cell(20,5).Range.Text = "48-46=2"
